# Auto-generated script applying numeric corrections to Leve profit calculations
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 8335390.5
$ws.Range("J17").Value = 8335390.5
$ws.Range("L17").Value = 25006171.5
$ws.Range("N17").Value = -25006507.5
$ws.Range("H40").Value = 3954.04
$ws.Range("I40").Value = 2383
$ws.Range("K40").Value = 2383
$ws.Range("M40").Value = -2208
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -931
$ws.Range("N43").ClearContents()
$ws.Range("H81").Value = 79987.5
$ws.Range("J81").Value = 79987.5
$ws.Range("L81").Value = 79987.5
$ws.Range("N81").Value = -81983.5
$ws.Range("H84").Value = 79987.5
$ws.Range("J84").Value = 79987.5
$ws.Range("L84").Value = 239962.5
$ws.Range("N84").Value = -249946.5
$ws.Range("H112").Value = 3197.0908
$ws.Range("I112").Value = 1989.3334
$ws.Range("J112").Value = 3650
$ws.Range("K112").Value = 5968.0002
$ws.Range("L112").Value = 10950
$ws.Range("M112").Value = -4860.0002
$ws.Range("N112").Value = -13166
$ws.Range("H138").Value = 3385.276
$ws.Range("J138").Value = 3799.5454
$ws.Range("L138").Value = 11398.6362
$ws.Range("N138").Value = -21678.6362

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5654.4614
$ws.Range("I2").Value = 5062.125
$ws.Range("K2").Value = 5062.125
$ws.Range("M2").Value = -4949.125
$ws.Range("H43").Value = 58788.2
$ws.Range("J43").Value = 54325
$ws.Range("L43").Value = 54325
$ws.Range("N43").Value = -54951
$ws.Range("H45").Value = 3990
$ws.Range("I45").Value = 3133.8572
$ws.Range("J45").Value = 5488.25
$ws.Range("K45").Value = 3133.8572
$ws.Range("L45").Value = 5488.25
$ws.Range("M45").Value = -2756.8572
$ws.Range("N45").Value = -6242.25
$ws.Range("H46").Value = 7059.8335
$ws.Range("J46").Value = 7094.7
$ws.Range("L46").Value = 7094.7
$ws.Range("N46").Value = -7732.7
$ws.Range("H61").Value = 8463.305
$ws.Range("I61").Value = 8463.305
$ws.Range("K61").Value = 8463.305
$ws.Range("M61").Value = -8251.305
$ws.Range("H74").Value = 5236.2666
$ws.Range("I74").Value = 2808.2173
$ws.Range("J74").Value = 13214.143
$ws.Range("K74").Value = 2808.2173
$ws.Range("L74").Value = 13214.143
$ws.Range("M74").Value = -1934.2173
$ws.Range("N74").Value = -14962.143
$ws.Range("H77").Value = 5236.2666
$ws.Range("I77").Value = 2808.2173
$ws.Range("J77").Value = 13214.143
$ws.Range("K77").Value = 14041.0865
$ws.Range("L77").Value = 66070.715
$ws.Range("M77").Value = -9673.086499999999
$ws.Range("N77").Value = -74806.715
$ws.Range("H116").Value = 5654.4614
$ws.Range("I116").Value = 5062.125
$ws.Range("K116").Value = 5062.125
$ws.Range("M116").Value = -2768.125
$ws.Range("H136").Value = 8463.305
$ws.Range("I136").Value = 8463.305
$ws.Range("K136").Value = 25389.915
$ws.Range("M136").Value = -22839.915

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5654.4614
$ws.Range("I3").Value = 5062.125
$ws.Range("K3").Value = 5062.125
$ws.Range("M3").Value = -4948.125
$ws.Range("H20").Value = 2906.7778
$ws.Range("I20").Value = 2734.625
$ws.Range("J20").Value = 3044.5
$ws.Range("K20").Value = 2734.625
$ws.Range("L20").Value = 3044.5
$ws.Range("M20").Value = -2487.625
$ws.Range("N20").Value = -3538.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3808.1904
$ws.Range("I16").Value = 3198.4707
$ws.Range("K16").Value = 3198.4707
$ws.Range("M16").Value = -2911.4707
$ws.Range("H58").Value = 10462.96
$ws.Range("I58").Value = 6116.1
$ws.Range("J58").Value = 13360.866
$ws.Range("K58").Value = 6116.1
$ws.Range("L58").Value = 13360.866
$ws.Range("M58").Value = -5913.1
$ws.Range("N58").Value = -13766.866
$ws.Range("H113").Value = 3808.1904
$ws.Range("I113").Value = 3198.4707
$ws.Range("K113").Value = 3198.4707
$ws.Range("M113").Value = -1028.4707
$ws.Range("H122").Value = 2675.889
$ws.Range("I122").Value = 2472.6191
$ws.Range("K122").Value = 7417.8573
$ws.Range("M122").Value = -4967.8573
$ws.Range("H132").Value = 4190.9287
$ws.Range("I132").Value = 2627.4
$ws.Range("J132").Value = 8099.75
$ws.Range("K132").Value = 7882.200000000001
$ws.Range("L132").Value = 24299.25
$ws.Range("M132").Value = -5352.200000000001
$ws.Range("N132").Value = -29359.25
$ws.Range("H134").Value = 5019.9473
$ws.Range("I134").Value = 2936.8462
$ws.Range("J134").Value = 9533.333000000001
$ws.Range("K134").Value = 8810.5386
$ws.Range("L134").Value = 28599.999
$ws.Range("M134").Value = -6275.5386
$ws.Range("N134").Value = -33669.999
$ws.Range("H136").Value = 10462.96
$ws.Range("I136").Value = 6116.1
$ws.Range("J136").Value = 13360.866
$ws.Range("K136").Value = 18348.3
$ws.Range("L136").Value = 40082.598
$ws.Range("M136").Value = -15798.3
$ws.Range("N136").Value = -45182.598

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13891130
$ws.Range("I131").Value = 45455650
$ws.Range("J131").Value = 2741.76
$ws.Range("K131").Value = 136366950
$ws.Range("L131").Value = 8225.280000000001
$ws.Range("M131").Value = -136361910
$ws.Range("N131").Value = -18305.28
$ws.Range("H137").Value = 5481.6665
$ws.Range("I137").Value = 3593.5
$ws.Range("K137").Value = 10780.5
$ws.Range("M137").Value = -5680.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15771
$ws.Range("J70").Value = 15771
$ws.Range("L70").Value = 15771
$ws.Range("N70").Value = -16311
$ws.Range("H73").Value = 15771
$ws.Range("J73").Value = 15771
$ws.Range("L73").Value = 15771
$ws.Range("N73").Value = -17643

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5496.3335
$ws.Range("I14").Value = 5495
$ws.Range("J14").Value = 5497
$ws.Range("K14").Value = 5495
$ws.Range("L14").Value = 5497
$ws.Range("M14").Value = -5327
$ws.Range("N14").Value = -5833
$ws.Range("H15").Value = 6991.2
$ws.Range("I15").Value = 7000
$ws.Range("K15").Value = 7000
$ws.Range("M15").Value = -6712
